# Insert a new data row at row 350 (pushing existing rows 350..437 down to 351..438)
# and populate it with the new record's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(350).Insert()

$ws.Cells.Item(350, 1).Value = 3
$ws.Cells.Item(350, 2).Value = 'Femacal de La Calera'
$ws.Cells.Item(350, 3).Value = 'Coquimbo'
$ws.Cells.Item(350, 4).Value = 44782
$ws.Cells.Item(350, 5).Value = 5
$ws.Cells.Item(350, 6).Value = 100112017
$ws.Cells.Item(350, 7).Value = 'Apio'
$ws.Cells.Item(350, 8).Value = 'Americana (o)'
$ws.Cells.Item(350, 9).Value = 'Primera'
$ws.Cells.Item(350, 10).Value = 230
$ws.Cells.Item(350, 11).Value = 9000
$ws.Cells.Item(350, 12).Value = 9500
$ws.Cells.Item(350, 13).Value = 9261
$ws.Cells.Item(350, 14).Value = '$/docena de matas'
$ws.Cells.Item(350, 15).Value = 'Provincia de Limarí'
$ws.Cells.Item(350, 16).Value = 1544
$ws.Cells.Item(350, 17).Value = 6
$ws.Cells.Item(350, 18).Value = 'Hortaliza'
